$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "protein_fasta" column (E) header/value is no longer used - clear it out.
# Clearing (rather than just removing the value) drops the cell's formatting too,
# so the now-empty E1 cell is omitted entirely when the sheet is saved - matching
# the rest of the (always empty) cells in column E.
$ws.Range("E1").Clear()

# Move the active selection to E1, as reflected in the saved view state.
$ws.Range("E1").Select()
